# Update the date header paragraph
$d = $word.ActiveDocument
$d.Paragraphs.Item(1).Range.Text = "2024-12-29 Sunday"

# Update each multiplication-table cell with its new value.
# Cells are addressed directly by (row, column) so that values which
# coincidentally match another cell's old/new text never get cross-matched.
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "44×52=2288"
$t.Cell(1, 2).Range.Text = "61×27=1647"
$t.Cell(1, 3).Range.Text = "12×97=1164"
$t.Cell(1, 4).Range.Text = "47×58=2726"
$t.Cell(1, 5).Range.Text = "89×62=5518"

$t.Cell(5, 1).Range.Text = "38×89=3382"
$t.Cell(5, 2).Range.Text = "43×74=3182"
$t.Cell(5, 3).Range.Text = "93×65=6045"
$t.Cell(5, 4).Range.Text = "40×61=2440"
$t.Cell(5, 5).Range.Text = "36×14=504"

$t.Cell(10, 1).Range.Text = "16×14=224"
$t.Cell(10, 2).Range.Text = "20×78=1560"
$t.Cell(10, 3).Range.Text = "78×34=2652"
$t.Cell(10, 4).Range.Text = "37×17=629"
$t.Cell(10, 5).Range.Text = "56×73=4088"

$t.Cell(15, 1).Range.Text = "41×48=1968"
$t.Cell(15, 2).Range.Text = "57×84=4788"
$t.Cell(15, 3).Range.Text = "35×86=3010"
$t.Cell(15, 4).Range.Text = "53×80=4240"
$t.Cell(15, 5).Range.Text = "14×35=490"

$t.Cell(20, 1).Range.Text = "79×56=4424"
$t.Cell(20, 2).Range.Text = "17×62=1054"
$t.Cell(20, 3).Range.Text = "31×20=620"
$t.Cell(20, 4).Range.Text = "64×45=2880"
$t.Cell(20, 5).Range.Text = "14×94=1316"
